$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header timestamp (row 1)
$ws.Cells.Item(1,1).Value = "Datos actualizados a 20 de Marzo de 2020 a las 13:46"

# Update rows whose province label and/or figures changed
$ws.Cells.Item(8,1).Value = "Navarra"
$ws.Cells.Item(8,2).Value = 593
$ws.Cells.Item(8,3).Value = 2
$ws.Cells.Item(8,4).Value = 585
$ws.Cells.Item(8,5).Value = 6
$ws.Cells.Item(21,1).Value = "Caceres"
$ws.Cells.Item(21,2).Value = 206
$ws.Cells.Item(21,3).Value = 2
$ws.Cells.Item(21,4).Value = 194
$ws.Cells.Item(21,5).Value = 10
$ws.Cells.Item(22,1).Value = "Guadalajara"
$ws.Cells.Item(22,2).Value = 205
$ws.Cells.Item(22,3).Value = 2
$ws.Cells.Item(22,4).Value = 200
$ws.Cells.Item(22,5).Value = 3
$ws.Cells.Item(23,1).Value = "Illes Balears"
$ws.Cells.Item(23,2).Value = 203
$ws.Cells.Item(23,3).Value = 10
$ws.Cells.Item(23,4).Value = 189
$ws.Cells.Item(23,5).Value = 4
$ws.Cells.Item(24,1).Value = "Pontevedra"
$ws.Cells.Item(24,2).Value = 193
$ws.Cells.Item(24,3).Value = 0
$ws.Cells.Item(24,4).Value = 191
$ws.Cells.Item(24,5).Value = 2
$ws.Cells.Item(25,1).Value = "Tenerife"
$ws.Cells.Item(25,2).Value = 192
$ws.Cells.Item(25,3).Value = 4
$ws.Cells.Item(25,4).Value = 185
$ws.Cells.Item(25,5).Value = 3
$ws.Cells.Item(42,1).Value = "Gran Canaria"
$ws.Cells.Item(42,2).Value = 70
$ws.Cells.Item(42,3).Value = 0
$ws.Cells.Item(42,4).Value = 69
$ws.Cells.Item(42,5).Value = 1
$ws.Cells.Item(43,1).Value = "Ourense"
$ws.Cells.Item(43,2).Value = 63
$ws.Cells.Item(43,3).Value = 0
$ws.Cells.Item(43,4).Value = 63
$ws.Cells.Item(43,5).Value = 0
$ws.Cells.Item(44,1).Value = "Avila"
$ws.Cells.Item(44,2).Value = 59
$ws.Cells.Item(44,3).Value = 2
$ws.Cells.Item(44,4).Value = 55
$ws.Cells.Item(44,5).Value = 2
$ws.Cells.Item(45,1).Value = "Igualada, Vilanova del Cami, Santa Margarida de Montbui y Odena"
$ws.Cells.Item(45,2).Value = 58
$ws.Cells.Item(45,3).Value = 0
$ws.Cells.Item(45,4).Value = 58
$ws.Cells.Item(45,5).Value = 3
$ws.Cells.Item(46,1).Value = "Soria"
$ws.Cells.Item(46,2).Value = 58
$ws.Cells.Item(46,3).Value = 4
$ws.Cells.Item(46,4).Value = 49
$ws.Cells.Item(46,5).Value = 5
$ws.Cells.Item(55,1).Value = "Fuerteventura"
$ws.Cells.Item(55,2).Value = 12
$ws.Cells.Item(55,3).Value = 0
$ws.Cells.Item(55,4).Value = 12
$ws.Cells.Item(55,5).Value = 0
$ws.Cells.Item(56,1).Value = "La Palma"
$ws.Cells.Item(56,2).Value = 7
$ws.Cells.Item(56,3).Value = 0
$ws.Cells.Item(56,4).Value = 7
$ws.Cells.Item(56,5).Value = 0
$ws.Cells.Item(57,1).Value = "Arroyo de la Luz"
$ws.Cells.Item(57,2).Value = 7
$ws.Cells.Item(57,3).Value = 0
$ws.Cells.Item(57,4).Value = 7
$ws.Cells.Item(57,5).Value = 0
$ws.Cells.Item(58,1).Value = "Ceuta"
$ws.Cells.Item(58,2).Value = 5
$ws.Cells.Item(58,3).Value = 0
$ws.Cells.Item(58,4).Value = 5
$ws.Cells.Item(58,5).Value = 0

# Append new row 61 for El Hierro
$ws.Cells.Item(61,1).Value = "El Hierro"
$ws.Cells.Item(61,2).Value = 1
$ws.Cells.Item(61,3).Value = 0
$ws.Cells.Item(61,4).Value = 1
$ws.Cells.Item(61,5).Value = 0
